$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row with the moved/renamed lightmap entry
$ws.Range("A14").Value = "Reimproved lighting (Vanilla, OptiFine)"

# Reflect the last active selection recorded in the saved file
$ws.Range("G13").Select()
